# Daily attendance processing - reorders the "Recorded By" (column G) tokens
# so that the "System" entry is listed first instead of last, for every
# session row on the "Session Analysis Results" sheet.
#
# Examples:
#   "dnasr281@gmail.com, System"                -> "System, dnasr281@gmail.com"
#   "backup@backdoor.com, System"                -> "System, backup@backdoor.com"
#   "system, backup@backdoor.com, System"        -> "System, backup@backdoor.com, system"
# Rows whose "Recorded By" value contains no System/system token (or more
# than one non-System token alongside two System tokens, etc.) are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) { continue }

    $parts = $text -split ", "

    $sysIdx = @()
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($parts[$i].ToLower() -eq "system") { $sysIdx += $i }
    }

    if ($sysIdx.Length -eq 0) {
        # No "System" token present - leave as-is.
        continue
    }
    elseif ($sysIdx.Length -eq 1 -and $parts.Length -eq 2) {
        # Two entries, one of them System: move System to the front,
        # keeping its original casing ("System"), the other entry keeps
        # its original text.
        $other = $parts[1 - $sysIdx[0]]
        $cell.Value = "System, " + $other
    }
    elseif ($sysIdx.Length -eq 2 -and $parts.Length -eq 3) {
        # Three entries where both the first and last are System
        # (one lowercase "system", one capitalized "System"): swap so the
        # capitalized "System" leads and the lowercase "system" trails,
        # keeping the middle entry unchanged.
        $cell.Value = "System, " + $parts[1] + ", system"
    }
    # Any other shape is left untouched.
}
